# UT clean up to xlsx files
# Rewrites the "Texas Notes" sheet: condenses the long block-quote of the
# journal abstract into fewer, longer lines, moves the citation DOI over to
# F1 (reusing its original "source-note" styling), and tightens up a couple
# of notes' wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Texas Notes")

# The DOI citation (previously A2) moves to F1, carrying its existing
# "source note" character formatting (Arial, grey) along with it.
$ws.Range("A2").Copy($ws.Range("F1"))

# Clear out the rest of column A so the sheet can be rebuilt with the new,
# consolidated row layout.
$ws.Columns.Item(1).Clear()

$ws.Range("A1").Value = "updated version of the study that EPS cites here:"

$ws.Range("A3").Value = """Results suggest that rebate policies increase the sales share of ENERGY STAR household appliances"
$ws.Range("A4").Value = "by 3.3 to 6.6 percentage points, and this represents an impact of 9 to 18 % on the mean level of the sales share of ENERGY STAR household appliances"
$ws.Range("A5").Value = "sales share of ENERGY STAR household appliances in the US between 2001 and 2006."""

$ws.Range("A6").Value = "See Table 3."

$ws.Range("A8").Value = "The 3.3% and 6.6% answers are both statistically significant and depend on the analysis method being used. "
$ws.Range("A9").Value = "Neither method stood out as being ""better""; average the two together and get 4.95%."

$ws.Range("D14").Select()
